$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 114
$ws.Range("H114").Value = 0.3

# Row 116
$ws.Range("G116").Value = 1.5

# Row 117
$ws.Range("H117").Value = 2.2

# Row 118
$ws.Range("G118").Value = 1.8
$ws.Range("H118").Value = 2.7
$ws.Range("V118").Value = 2.4

# Row 119
$ws.Range("G119").Value = 1.5
$ws.Range("H119").Value = 0.8
$ws.Range("V119").Value = 1.2

# Row 120
$ws.Range("G120").Value = 1.7
$ws.Range("H120").Value = 1.1
$ws.Range("V120").Value = 3.3

# Row 121
$ws.Range("G121").Value = 1.1
$ws.Range("H121").Value = -0.1
$ws.Range("V121").Value = 2

# Row 122
$ws.Range("G122").Value = -3
$ws.Range("H122").Value = -13.2
$ws.Range("V122").Value = -3.6

# Row 123
$ws.Range("F123").Value = -28.1
$ws.Range("G123").Value = -14.5
$ws.Range("H123").Value = -39.1
$ws.Range("V123").Value = -29.8

# Row 124
$ws.Range("E124").Value = -5.5
$ws.Range("F124").Value = 23.2
$ws.Range("H124").Value = 60.9

# Row 125
$ws.Range("E125").Value = -0.9
$ws.Range("G125").Value = -4.4
$ws.Range("H125").Value = -1.7
$ws.Range("V125").Value = -1.8

# Row 126
$ws.Range("F126").Value = -4.2
$ws.Range("G126").Value = -1.2
$ws.Range("H126").Value = -1.1
$ws.Range("V126").Value = 4.5

# Row 127
$ws.Range("E127").Value = 7.6
$ws.Range("F127").Value = 1.9
$ws.Range("G127").Value = 14.3
$ws.Range("H127").Value = 9.2
$ws.Range("N127").Value = 10.5
$ws.Range("V127").Value = 41.9
